$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace "Aud S" (and one "Aud J") with "Aud Q" for the five sessions that moved rooms.
$ws.Range("D4").Value = "29.01: <strong>Collaborative learning session 2</strong> in Aud Q"
$ws.Range("D5").Value = "05.02: <strong>Exercise session 2</strong> in Aud Q"
$ws.Range("D6").Value = "12.02: <strong>Case session 2</strong> in Aud Q"
$ws.Range("D7").Value = "19.02: <strong>Collaborative learning session 3</strong> in Aud Q"
$ws.Range("D8").Value = "26.02: <strong>Case session 3</strong>  in Aud Q"

# Update the saved view state to match: D9 selected (sheet was scrolled back to top).
$ws.Range("D9").Select()
